# Apply the edits described in the diff:
#
# 1. Sheet "20_Properties of Circles" (1st sheet): the numbers in A4:A10 were
#    mis-keyed as 4,5,6,7,8,9,0 and are corrected to 3,4,5,6,7,8,9.
# 2. Sheet "21_Prop of Tangent to Circle" (2nd sheet): same renumbering fix
#    in A4:A10.
# 3. The view switches from sheet 1 being the selected/active tab to sheet 2
#    being the selected/active tab, and the selections on both sheets grow
#    from a single cell (A2) to the full data column (A2:A17 / A2:A14).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("20_Properties of Circles")
$ws2 = $wb.Worksheets.Item("21_Prop of Tangent to Circle")

# --- Fix up the numbering in column A, rows 4-10, on both sheets ---
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("A4").Value = 3
    $ws.Range("A5").Value = 4
    $ws.Range("A6").Value = 5
    $ws.Range("A7").Value = 6
    $ws.Range("A8").Value = 7
    $ws.Range("A9").Value = 8
    $ws.Range("A10").Value = 9
}

# --- Update the selection shown on sheet 1 (no longer the active tab) ---
$ws1.Activate()
$ws1.Range("A2:A17").Select()

# --- Update the selection on sheet 2 and make it the active/selected tab ---
$ws2.Activate()
$ws2.Range("A2:A14").Select()

$wb.Save()
